$d = $word.ActiveDocument
$wdNS = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Replace-WholeParagraph($anchorText, $newParaInnerXml) {
    # Find the paragraph that contains $anchorText, insert a brand-new
    # paragraph (built from raw OOXML) immediately in front of it, then
    # delete the original paragraph (including its paragraph mark).
    $found = $d.Content
    $found.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $para = $found.Paragraphs(1)
    $pr = $para.Range

    $insPoint = $d.Range($pr.Start, $pr.Start)
    $insPoint.InsertXML("<w:p $wdNS>$newParaInnerXml</w:p>")

    $found2 = $d.Content
    $found2.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $para2 = $found2.Paragraphs(1)
    $para2.Range.Delete()
}

# --- 1) "Código" header cell: vertically center the cell content ---
$t1 = $d.Tables(1)
$t1.Cell(1, 1).VerticalAlignment = 1

# --- 2) Header cell "Requerimiento funcional" -> single bold run "Requerimiento Funcional " ---
Replace-WholeParagraph "Requerimiento funcional" '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Requerimiento Funcional </w:t></w:r>'

# --- 3) "Usuarios: " + "cliente,administradores" + ", jugadores" -> restructured run/proofErr layout ---
Replace-WholeParagraph "Usuarios: cliente,administradores, jugadores" '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>Usuarios:profesores</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>,administradores</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, jugadores</w:t></w:r>'

# --- 4) "Usuarios:cliente" (no space, proofErr-wrapped) -> "Usuarios:profesores" ---
$d.Content.Find.Execute("Usuarios:cliente", $true, $false, $false, $false, $false, $true, 1, $false, "Usuarios:profesores", 2)

# --- 5) "Usuarios: cliente, administradores" (3x) -> "Usuarios: profesores, administradores" ---
$d.Content.Find.Execute("Usuarios: cliente, administradores", $true, $false, $false, $false, $false, $true, 1, $false, "Usuarios: profesores, administradores", 2)

# --- 6) remaining standalone "Usuarios: cliente" -> "Usuarios: profesores" ---
$d.Content.Find.Execute("Usuarios: cliente", $true, $false, $false, $false, $false, $true, 1, $false, "Usuarios: profesores", 2)

# --- 7) "Usuarios: Clientes, Administradores" -> "Usuarios: profesores, Administradores, Jugadores" ---
$d.Content.Find.Execute("Usuarios: Clientes, Administradores", $true, $false, $false, $false, $false, $true, 1, $false, "Usuarios: profesores, Administradores, Jugadores", 2)

# --- 8) final "Usuarios: Clientes" cell -> new text + first-line indent paragraph property ---
Replace-WholeParagraph "Usuarios: Clientes" '<w:pPr><w:ind w:firstLine="708"/></w:pPr><w:r><w:t>Usuarios: profesores, Administradores, Jugadores</w:t></w:r>'
